# homeword excel sheet has been modified
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row update (shorten the vulnerability description) ---
$ws.Range("B1").Value = "potential vulnerablity (Violates user privacy/Causes security threat)"

# --- Column C (Manual / Automatic) filled top-to-bottom first ---
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 3).Value = "Manual"
}
for ($r = 22; $r -le 41; $r++) {
    $ws.Cells.Item($r, 3).Value = "Automatic"
}

# --- Column B (Violates User Privacy / Cause Security Threat) filled next ---
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = "Violates User Privacy"
}
for ($r = 12; $r -le 21; $r++) {
    $ws.Cells.Item($r, 2).Value = "Cause Security Threat"
}
for ($r = 22; $r -le 31; $r++) {
    $ws.Cells.Item($r, 2).Value = "Violates User Privacy"
}
for ($r = 32; $r -le 41; $r++) {
    $ws.Cells.Item($r, 2).Value = "Cause Security Threat"
}

# --- New rows 32-41 need a Chain# in column A as well ---
for ($r = 32; $r -le 41; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# --- Column width adjustments ---
# (target widths of 20.5 / 20.33203125 / 19.33203125 / 18.83203125 are reproduced
# as closely as the ColumnWidth COM property's internal pixel quantization allows)
$ws.Columns.Item(7).ColumnWidth = 19.666666666666668
$ws.Columns.Item(8).ColumnWidth = 19.5
$ws.Columns.Item(9).ColumnWidth = 18.5
$ws.Columns.Item(11).ColumnWidth = 18.0
$ws.Columns.Item(13).ColumnWidth = 19.666666666666668

# --- Selection change ---
$ws.Range("C22").Select()
